$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("A2").Value = "B"
$ws.Range("A3").Value = "C"
$ws.Range("A4").Value = "D"
$ws.Range("A5").Value = "E"

$ws.Range("A6").Select() | Out-Null
